$wb = $excel.ActiveWorkbook

# --- "optimize" sheet: collapse the iteration-count parameters down to a single
#     placeholder value. These are text cells (shared strings), so an apostrophe
#     prefix is used to keep "1000" stored as text instead of a number, and the
#     formatting that the text-coercion leaves behind is cleared right after. ---
$optimize = $wb.Worksheets.Item("optimize")
$optimize.Range("B7").Formula = "'1000"
$optimize.Range("B7").ClearFormats()
$optimize.Range("B8").Formula = "'1000"
$optimize.Range("B8").ClearFormats()

# --- "studio" sheet: point the saved tool/view state at analysis instead of optimize ---
$studio = $wb.Worksheets.Item("studio")
$studio.Range("B2").Value = "analysis"
$studio.Range("B3").Value = "Reconstructed"

# --- Move the "landmarks" sheet to the end of the workbook (recreated, so it
#     picks up a fresh sheetId, matching how Excel reassigns one when a tab is
#     removed and re-added) ---
$landmarks = $wb.Worksheets.Item("landmarks")
$domain = $landmarks.Range("A1").Value()
$name = $landmarks.Range("B1").Value()
$visible = $landmarks.Range("C1").Value()
$color = $landmarks.Range("D1").Value()
$comment = $landmarks.Range("E1").Value()
$landmarks.Delete()

$lastIndex = $wb.Worksheets.Count
$newLandmarks = $wb.Worksheets.Add($null, $wb.Worksheets.Item($lastIndex))
$newLandmarks.Name = "landmarks"
$newLandmarks.Range("A1").Value = $domain
$newLandmarks.Range("B1").Value = $name
$newLandmarks.Range("C1").Value = $visible
$newLandmarks.Range("D1").Value = $color
$newLandmarks.Range("E1").Value = $comment
